$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1604938271604938
$ws.Range("C2").Value = 0.6296296296296297
$ws.Range("J2").Value = 0.01646090534979424
$ws.Range("P2").Value = 0.1275720164609054
$ws.Range("S2").Value = 0.06584362139917696
$ws.Range("B3").Value = 0.01257861635220126
$ws.Range("C3").Value = 0.0440251572327044
$ws.Range("J3").Value = 0.03144654088050314
$ws.Range("P3").Value = 0.7672955974842768
$ws.Range("S3").Value = 0.1446540880503145
$ws.Range("J4").Value = 0.1052631578947368
$ws.Range("P4").Value = 0.6842105263157895
$ws.Range("S4").Value = 0.2105263157894737
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.25
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.04878048780487805
$ws.Range("F6").Value = 0.05365853658536585
$ws.Range("J6").Value = 0.2292682926829268
$ws.Range("O6").Value = 0.00975609756097561
$ws.Range("Q6").Value = 0.1121951219512195
$ws.Range("R6").Value = 0.08780487804878048
$ws.Range("S6").Value = 0.4585365853658537
$ws.Range("B7").Value = 0.1186440677966102
$ws.Range("D7").Value = 0.02542372881355932
$ws.Range("F7").Value = 0.01694915254237288
$ws.Range("J7").Value = 0.1694915254237288
$ws.Range("O7").Value = 0.0423728813559322
$ws.Range("Q7").Value = 0.1101694915254237
$ws.Range("R7").Value = 0.05932203389830509
$ws.Range("S7").Value = 0.4576271186440678
$ws.Range("B8").Value = 0.05908096280087528
$ws.Range("D8").Value = 0.01094091903719912
$ws.Range("E8").Value = 0.002188183807439825
$ws.Range("F8").Value = 0.07439824945295405
$ws.Range("J8").Value = 0.1487964989059081
$ws.Range("O8").Value = 0.0350109409190372
$ws.Range("Q8").Value = 0.1509846827133479
$ws.Range("R8").Value = 0.1115973741794311
$ws.Range("S8").Value = 0.4070021881838075
$ws.Range("B9").Value = 0.1017543859649123
$ws.Range("D9").Value = 0.003508771929824561
$ws.Range("E9").Value = 0.003508771929824561
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1403508771929824
$ws.Range("O9").Value = 0.04210526315789474
$ws.Range("Q9").Value = 0.1684210526315789
$ws.Range("R9").Value = 0.09473684210526316
$ws.Range("S9").Value = 0.3789473684210526
$ws.Range("B10").Value = 0.0933852140077821
$ws.Range("D10").Value = 0.02334630350194553
$ws.Range("E10").Value = 0.001556420233463035
$ws.Range("F10").Value = 0.0622568093385214
$ws.Range("J10").Value = 0.1066147859922179
$ws.Range("O10").Value = 0.02334630350194553
$ws.Range("Q10").Value = 0.1961089494163424
$ws.Range("R10").Value = 0.1035019455252918
$ws.Range("S10").Value = 0.3898832684824903
$ws.Range("G11").Value = 0.08928571428571429
$ws.Range("J11").Value = 0.119047619047619
$ws.Range("K11").Value = 0.119047619047619
$ws.Range("L11").Value = 0.6726190476190477
$ws.Range("G12").Value = 0.7203389830508474
$ws.Range("J12").Value = 0.2203389830508475
$ws.Range("K12").Value = 0.01694915254237288
$ws.Range("L12").Value = 0.0423728813559322
$ws.Range("G13").Value = 0.5405405405405406
$ws.Range("J13").Value = 0.4324324324324325
$ws.Range("S13").Value = 0.02702702702702703
$ws.Range("F15").Value = 0.003968253968253968
$ws.Range("H15").Value = 0.1626984126984127
$ws.Range("I15").Value = 0.0873015873015873
$ws.Range("J15").Value = 0.3452380952380952
$ws.Range("K15").Value = 0.04365079365079365
$ws.Range("M15").Value = 0.01587301587301587
$ws.Range("O15").Value = 0.03968253968253968
$ws.Range("S15").Value = 0.3015873015873016
$ws.Range("F16").Value = 0.01149425287356322
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("I16").Value = 0.1436781609195402
$ws.Range("J16").Value = 0.4080459770114943
$ws.Range("K16").Value = 0.08620689655172414
$ws.Range("O16").Value = 0.05172413793103448
$ws.Range("S16").Value = 0.132183908045977
$ws.Range("F17").Value = 0.01990049751243781
$ws.Range("H17").Value = 0.1691542288557214
$ws.Range("I17").Value = 0.1194029850746269
$ws.Range("J17").Value = 0.4502487562189055
$ws.Range("K17").Value = 0.06716417910447761
$ws.Range("M17").Value = 0.01492537313432836
$ws.Range("O17").Value = 0.06467661691542288
$ws.Range("S17").Value = 0.0945273631840796
$ws.Range("F18").Value = 0.02991452991452992
$ws.Range("H18").Value = 0.1623931623931624
$ws.Range("I18").Value = 0.1282051282051282
$ws.Range("J18").Value = 0.4230769230769231
$ws.Range("K18").Value = 0.05982905982905983
$ws.Range("M18").Value = 0.008547008547008548
$ws.Range("N18").Value = 0.004273504273504274
$ws.Range("O18").Value = 0.07264957264957266
$ws.Range("S18").Value = 0.1111111111111111
$ws.Range("F19").Value = 0.0126883425852498
$ws.Range("H19").Value = 0.2268041237113402
$ws.Range("I19").Value = 0.1292624900872324
$ws.Range("J19").Value = 0.3798572561459159
$ws.Range("K19").Value = 0.06026962727993656
$ws.Range("M19").Value = 0.01982553528945281
$ws.Range("N19").Value = 0.0007930214115781126
$ws.Range("O19").Value = 0.0761300555114988
$ws.Range("S19").Value = 0.0943695479777954
